$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old "Rank (code churn)" formulas ranked every row (2-13) against the
# whole C2:C13 column. They are replaced by a narrower analysis scoped to
# rows 8-12 (paired with a new "Rank (branch coverage)" column), so remove
# the stale formulas from rows 2-7 and 13.
$ws.Range("E2:E7").ClearContents()
$ws.Range("E13").ClearContents()

# The lone Overall Branch Coverage value that used to sit on row 13 moves
# up into the rows 8-12 block below.
$ws.Range("B13").ClearContents()

# Overall Branch Coverage values for rows 8-12
$ws.Range("B8").Value = 0.77
$ws.Range("B9").Value = 0.78
$ws.Range("B10").Value = 0.81
$ws.Range("B11").Value = 0.81
$ws.Range("B12").Value = 0.82

# Rank (branch coverage), column D, rows 8-12
$ws.Range("D8").Formula = "=_xlfn.RANK.AVG(B8,`$B`$8:`$B`$12,1)"
$ws.Range("D9:D12").Formula = "=_xlfn.RANK.AVG(B9,`$B`$8:`$B`$12,1)"

# Rank (code churn), column E, rows 8-12, rescoped to C8:C12
$ws.Range("E8").Formula = "=_xlfn.RANK.AVG(C8,C8:C12,1)"
$ws.Range("E9").Formula = "=_xlfn.RANK.AVG(C9,`$C`$8:`$C`$12,1)"

# Square of difference, column F, rows 8-12
$ws.Range("F8").Formula = "=(D8-E8)^2"
$ws.Range("F9:F12").Formula = "=(D9-E9)^2"

$ws.Range("E10:E12").Formula = "=_xlfn.RANK.AVG(C10,`$C`$8:`$C`$12,1)"

# Sum of squared differences, and the resulting Spearman correlation
$ws.Range("F13").Formula = "=SUM(F8:F12)"
$ws.Range("F14").Formula = "=1-((6*F13)/(125-5))"

$ws.Range("F15").Select() | Out-Null
